$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B: all values become 0
$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 0
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = 0
$ws.Range("B7").Value = 0
$ws.Range("B8").Value = 0
$ws.Range("B9").Value = 0

# Column C
$ws.Range("C2").Value = 0
$ws.Range("C3").Value = 0.6911885779249263
$ws.Range("C4").Value = 0
$ws.Range("C5").Value = 0.6862955558970697
$ws.Range("C6").Value = 0
$ws.Range("C7").Value = -0.7682308209596913
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = -0.7383293258771705

# Column D
$ws.Range("D2").Value = 0.6607615832614187
$ws.Range("D3").Value = -0.6265692778103581
$ws.Range("D4").Value = -0.6854857790183778
$ws.Range("D5").Value = -0.7070270614751112
$ws.Range("D6").Value = 0.6619919067668416
$ws.Range("D7").Value = 0.8091045540452172
$ws.Range("D8").Value = -0.7146818159952912
$ws.Range("D9").Value = 0.7562226113302229
